$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 29 data (appended at the end of the Artfynd sheet)
$ws.Range("A29").Value = 131114362
$ws.Range("B29").Value = 79243
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("P29").Value = "Lobergshöjden, Lobergshöjden, Vrm"
$ws.Range("Q29").Value = 477396
$ws.Range("R29").Value = 6591981
$ws.Range("S29").Value = 10
$ws.Range("T29").Value = "Örebro"
$ws.Range("U29").Value = "Karlskoga"
$ws.Range("V29").Value = "Värmland"
$ws.Range("W29").Value = "Karlskoga"

# Date/time-like values must stay as plain text, not be converted to date
# serials. Force a text number format before assigning, then clear the
# formatting again so no extra style is left behind on the cell.
$ws.Range("Y29").NumberFormat = "@"
$ws.Range("Y29").Value = "2026-02-11"
$ws.Range("Y29").ClearFormats()

$ws.Range("Z29").NumberFormat = "@"
$ws.Range("Z29").Value = "11:23"
$ws.Range("Z29").ClearFormats()

$ws.Range("AA29").NumberFormat = "@"
$ws.Range("AA29").Value = "2026-02-11"
$ws.Range("AA29").ClearFormats()

$ws.Range("AB29").NumberFormat = "@"
$ws.Range("AB29").Value = "11:23"
$ws.Range("AB29").ClearFormats()

$ws.Range("AD29").Value = $false
$ws.Range("AE29").Value = $false
$ws.Range("AG29").Value = $false

$ws.Range("AW29").Value = "Jim Hellquist"
$ws.Range("AX29").Value = "Jim Hellquist"
